$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.635.63'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '1.617.02'
$ws.Range("E3").Value = '  -1.16%  '
$ws.Range("E4").Value = '  -0.89%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.515'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.56%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.12'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("E9").Value = '  -1.26%  '
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("E11").Value = '  -1.13%  '
$ws.Range("D12").Value = '1.847.01'
$ws.Range("E12").Value = '  -1.14%  '
$ws.Range("D13").Value = '1.624.45'
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.555'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.58'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.30%  '
$ws.Range("D17").Value = '27.650.45'
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '226.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").Value = '0.0₃0714'
$ws.Range("E20").Value = '  -1.36%  '
$ws.Range("E21").Value = '  -0.93%  '
$ws.Range("E22").Value = '  -1.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.19%  '
$ws.Range("E24").Value = '  -1.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.991'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("E30").Value = '  -1.00%  '
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("D34").Value = '1.391.66'
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("E35").Value = '  +1.21%  '
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("E37").Value = '  -1.82%  '
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("E39").Value = '  -1.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.841'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.53%  '
$ws.Range("E41").Value = '  -1.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.990'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.95%  '
$ws.Range("E43").Value = '  -1.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.06%  '
$ws.Range("D46").Value = '1.756.46'
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("E47").Value = '  -4.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.57'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.92%  '
